# Update symbol list: refresh Price/Volume(1h)/Hora for rows 2-51,
# and correct Coin/Link for rows that were re-ranked (8-9, 19-24).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; D="325.63"; E="-0.92%"; G="2"}
    @{Row=3; D="43.93"; E="-2.03%"; G="2"}
    @{Row=4; D="5.502"; E="-1.90%"; G="2"}
    @{Row=5; D="0.07984"; E="-2.00%"; G="2"}
    @{Row=6; D="1.993"; E="4.20%"; G="2"}
    @{Row=7; D="4.301"; E="-0.55%"; G="2"}
    @{Row=8; B="MXToken"; C="https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; D="0.9478"; E="-0.68%"; G="2"}
    @{Row=9; B="BTSEToken"; C="https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"; D="2.554"; E="-7.30%"; G="2"}
    @{Row=10; D="0.1141"; E="-2.67%"; G="2"}
    @{Row=11; D="0.1831"; E="-4.35%"; G="2"}
    @{Row=12; D="11.79"; E="35.31%"; G="2"}
    @{Row=13; D="0.09609"; E="-3.37%"; G="2"}
    @{Row=14; D="0.04751"; E="13.68%"; G="2"}
    @{Row=15; D="0.1067"; E="0.20%"; G="2"}
    @{Row=16; D="0.001276"; E="0.38%"; G="2"}
    @{Row=17; D="0.04056"; E="-5.20%"; G="2"}
    @{Row=18; D="0.005752"; E="-5.65%"; G="2"}
    @{Row=19; B="LEO"; C="https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; D="3.367"; E="-6.11%"; G="2"}
    @{Row=20; B="BitpandaEcosystemToken"; C="https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"; D="0.3478"; E="-0.17%"; G="2"}
    @{Row=21; B="ProBitToken"; C="https://coinranking.com/coin/lQP4d6T2+probittoken-prob"; D="0.1410"; E="3.16%"; G="2"}
    @{Row=22; B="ZBToken"; C="https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"; D="0.2547"; E="-1.94%"; G="2"}
    @{Row=23; B="BitKan"; C="https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"; D="0.001243"; E="0.18%"; G="2"}
    @{Row=24; B="HotbitToken"; C="https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"; D="0.004294"; E="-6.79%"; G="2"}
    @{Row=25; E="-3.65%"; G="2"}
    @{Row=26; E="-6.45%"; G="2"}
    @{Row=27; G="2"}
    @{Row=28; G="2"}
    @{Row=29; G="2"}
    @{Row=30; G="2"}
    @{Row=31; G="2"}
    @{Row=32; G="2"}
    @{Row=33; G="2"}
    @{Row=34; G="2"}
    @{Row=35; G="2"}
    @{Row=36; G="2"}
    @{Row=37; G="2"}
    @{Row=38; D="0.02515"; E="-7.07%"; G="2"}
    @{Row=39; D="0.05529"; E="-1.88%"; G="2"}
    @{Row=40; D="0.007529"; E="-1.89%"; G="2"}
    @{Row=41; E="-0.87%"; G="2"}
    @{Row=42; D="0.007437"; E="-34.38%"; G="2"}
    @{Row=43; D="0.002014"; E="-3.73%"; G="2"}
    @{Row=44; D="0.008372"; E="-3.63%"; G="2"}
    @{Row=45; D="0.00007111"; E="0.15%"; G="2"}
    @{Row=46; E="-0.40%"; G="2"}
    @{Row=47; E="1.03%"; G="2"}
    @{Row=48; D="0.003527"; E="2.33%"; G="2"}
    @{Row=49; D="0.00002099"; E="-0.40%"; G="2"}
    @{Row=50; D="0.0001999"; E="-0.40%"; G="2"}
    @{Row=51; G="2"}
)

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) { $ws.Cells.Item($u.Row, 2).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Cells.Item($u.Row, 3).Value = $u.C }
    if ($u.ContainsKey("D")) { $ws.Cells.Item($u.Row, 4).Value = "'" + $u.D }
    if ($u.ContainsKey("E")) { $ws.Cells.Item($u.Row, 5).Value = "'" + $u.E }
    if ($u.ContainsKey("G")) { $ws.Cells.Item($u.Row, 7).Value = "'" + $u.G }
}
